$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns to match
# the latest scrape. D-column values that look like plain numbers are
# written with a leading apostrophe so Excel keeps them as text (matching
# the source data which stores these as inline strings, preserving exact
# formatting such as trailing zeros).

$ws.Cells.Item(2, 4).Value = "27.960.12"
$ws.Cells.Item(2, 5).Value = "  -2.09%  "
$ws.Cells.Item(3, 4).Value = "1.884.68"
$ws.Cells.Item(3, 5).Value = "  -1.42%  "
$ws.Cells.Item(4, 5).Value = "  +0.16%  "
$ws.Cells.Item(5, 4).Value = "'312.87"
$ws.Cells.Item(5, 5).Value = "  -0.68%  "
$ws.Cells.Item(6, 5).Value = "  +0.14%  "
$ws.Cells.Item(7, 4).Value = "'0.4997"
$ws.Cells.Item(7, 5).Value = "  -4.74%  "
$ws.Cells.Item(8, 4).Value = "'0.3850"
$ws.Cells.Item(8, 5).Value = "  -2.79%  "
$ws.Cells.Item(9, 4).Value = "'0.09141"
$ws.Cells.Item(9, 5).Value = "  -6.06%  "
$ws.Cells.Item(10, 4).Value = "'1.121"
$ws.Cells.Item(10, 5).Value = "  -2.70%  "
$ws.Cells.Item(11, 4).Value = "'41.67"
$ws.Cells.Item(11, 5).Value = "  -0.87%  "
$ws.Cells.Item(12, 4).Value = "'6.328"
$ws.Cells.Item(12, 5).Value = "  -3.32%  "
$ws.Cells.Item(13, 5).Value = "  -2.43%  "
$ws.Cells.Item(14, 4).Value = "1.880.19"
$ws.Cells.Item(14, 5).Value = "  -2.01%  "
$ws.Cells.Item(15, 4).Value = "'7.268"
$ws.Cells.Item(15, 5).Value = "  -3.93%  "
$ws.Cells.Item(16, 4).Value = "'1.002"
$ws.Cells.Item(16, 5).Value = "  +0.14%  "
$ws.Cells.Item(17, 4).Value = "'0.00001105"
$ws.Cells.Item(17, 5).Value = "  -3.03%  "
$ws.Cells.Item(18, 4).Value = "'91.38"
$ws.Cells.Item(18, 5).Value = "  -3.55%  "
$ws.Cells.Item(19, 4).Value = "'0.06627"
$ws.Cells.Item(19, 5).Value = "  -0.48%  "
$ws.Cells.Item(20, 4).Value = "'17.84"
$ws.Cells.Item(20, 5).Value = "  -2.09%  "
$ws.Cells.Item(21, 5).Value = "  +0.03%  "
$ws.Cells.Item(22, 4).Value = "'6.186"
$ws.Cells.Item(22, 5).Value = "  -2.27%  "
$ws.Cells.Item(23, 4).Value = "28.026.36"
$ws.Cells.Item(23, 5).Value = "  -2.19%  "
$ws.Cells.Item(24, 4).Value = "'11.34"
$ws.Cells.Item(24, 5).Value = "  -1.70%  "
$ws.Cells.Item(25, 4).Value = "'2.320"
$ws.Cells.Item(25, 5).Value = "  +0.72%  "
$ws.Cells.Item(26, 4).Value = "2.097.85"
$ws.Cells.Item(26, 5).Value = "  -1.78%  "
$ws.Cells.Item(27, 4).Value = "'2.540"
$ws.Cells.Item(27, 5).Value = "  -5.42%  "
$ws.Cells.Item(28, 4).Value = "'158.08"
$ws.Cells.Item(28, 5).Value = "  -0.43%  "
$ws.Cells.Item(29, 4).Value = "'20.71"
$ws.Cells.Item(29, 5).Value = "  -2.62%  "
$ws.Cells.Item(30, 4).Value = "'126.36"
$ws.Cells.Item(30, 5).Value = "  -2.11%  "
$ws.Cells.Item(31, 4).Value = "'1.065"
$ws.Cells.Item(31, 5).Value = "  -4.37%  "
$ws.Cells.Item(32, 4).Value = "'0.1052"
$ws.Cells.Item(32, 5).Value = "  -3.30%  "
$ws.Cells.Item(33, 4).Value = "'5.572"
$ws.Cells.Item(33, 5).Value = "  -3.32%  "
$ws.Cells.Item(34, 4).Value = "'3.606"
$ws.Cells.Item(34, 5).Value = "  -0.86%  "
$ws.Cells.Item(35, 4).Value = "'9.373"
$ws.Cells.Item(35, 5).Value = "  -5.76%  "
$ws.Cells.Item(36, 4).Value = "'0.06536"
$ws.Cells.Item(36, 5).Value = "  -3.54%  "
$ws.Cells.Item(37, 4).Value = "'0.02393"
$ws.Cells.Item(37, 5).Value = "  -1.78%  "
$ws.Cells.Item(38, 4).Value = "'1.302"
$ws.Cells.Item(38, 5).Value = "  +9.18%  "
$ws.Cells.Item(39, 4).Value = "'0.2179"
$ws.Cells.Item(39, 5).Value = "  -2.51%  "
$ws.Cells.Item(40, 4).Value = "'1.204"
$ws.Cells.Item(40, 5).Value = "  -4.99%  "
$ws.Cells.Item(41, 4).Value = "'0.6403"
$ws.Cells.Item(41, 5).Value = "  -1.28%  "
$ws.Cells.Item(42, 4).Value = "'11.50"
$ws.Cells.Item(42, 5).Value = "  -2.82%  "
$ws.Cells.Item(43, 4).Value = "'4.934"
$ws.Cells.Item(43, 5).Value = "  -3.05%  "
$ws.Cells.Item(44, 5).Value = "  +0.17%  "
$ws.Cells.Item(45, 4).Value = "'13.31"
$ws.Cells.Item(45, 5).Value = "  -2.14%  "
$ws.Cells.Item(46, 4).Value = "'0.6009"
$ws.Cells.Item(46, 5).Value = "  -1.61%  "
$ws.Cells.Item(47, 5).Value = "  +1.02%  "
$ws.Cells.Item(48, 4).Value = "'3.673"
$ws.Cells.Item(48, 5).Value = "  -2.17%  "
$ws.Cells.Item(49, 4).Value = "'1.990"
$ws.Cells.Item(49, 5).Value = "  -2.32%  "
$ws.Cells.Item(50, 4).Value = "'1.204"
$ws.Cells.Item(50, 5).Value = "  -0.43%  "
$ws.Cells.Item(51, 4).Value = "'121.01"
$ws.Cells.Item(51, 5).Value = "  -3.62%  "
